$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2: A=ECs D=FAPs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Rtn4"
$ws.Cells.Item(2,3).Value = "Tnfrsf19"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 22.91616466666667
$ws.Cells.Item(2,8).Value = 68.74849400000001
$ws.Cells.Item(2,9).Value = 0.05560000359891543
$ws.Cells.Item(2,10).Value = 0.05560000359891545
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.11008
$ws.Cells.Item(2,14).Value = 0.33024
$ws.Cells.Item(2,15).Value = 0.06903792193848049
$ws.Cells.Item(2,16).Value = 0.06903792193848049
$ws.Cells.Item(2,17).Value = 2.522611406506667
$ws.Cells.Item(2,18).Value = 22.70350265856
$ws.Cells.Item(2,19).Value = 0.003838508708241158
$ws.Cells.Item(2,20).Value = 0.003838508708241159

# row 3: A=ECs D=sCs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Rtn4"
$ws.Cells.Item(3,3).Value = "Tnfrsf19"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 22.91616466666667
$ws.Cells.Item(3,8).Value = 68.74849400000001
$ws.Cells.Item(3,9).Value = 0.05560000359891543
$ws.Cells.Item(3,10).Value = 0.05560000359891545
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.484406
$ws.Cells.Item(3,14).Value = 4.453218
$ws.Cells.Item(3,15).Value = 0.9309620780615195
$ws.Cells.Item(3,16).Value = 0.9309620780615195
$ws.Cells.Item(3,17).Value = 34.016892328188
$ws.Cells.Item(3,18).Value = 306.152030953692
$ws.Cells.Item(3,19).Value = 0.05176149489067428
$ws.Cells.Item(3,20).Value = 0.05176149489067429

# row 4: A=FAPs D=FAPs
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Rtn4"
$ws.Cells.Item(4,3).Value = "Tnfrsf19"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 93.34790299999999
$ws.Cells.Item(4,8).Value = 280.043709
$ws.Cells.Item(4,9).Value = 0.2264839609178002
$ws.Cells.Item(4,10).Value = 0.2264839609178003
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.11008
$ws.Cells.Item(4,14).Value = 0.33024
$ws.Cells.Item(4,15).Value = 0.06903792193848049
$ws.Cells.Item(4,16).Value = 0.06903792193848049
$ws.Cells.Item(4,17).Value = 10.27573716224
$ws.Cells.Item(4,18).Value = 92.48163446015998
$ws.Cells.Item(4,19).Value = 0.01563598201416096
$ws.Cells.Item(4,20).Value = 0.01563598201416096

# row 5: A=FAPs D=sCs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Rtn4"
$ws.Cells.Item(5,3).Value = "Tnfrsf19"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 93.34790299999999
$ws.Cells.Item(5,8).Value = 280.043709
$ws.Cells.Item(5,9).Value = 0.2264839609178002
$ws.Cells.Item(5,10).Value = 0.2264839609178003
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.484406
$ws.Cells.Item(5,14).Value = 4.453218
$ws.Cells.Item(5,15).Value = 0.9309620780615195
$ws.Cells.Item(5,16).Value = 0.9309620780615195
$ws.Cells.Item(5,17).Value = 138.566187300618
$ws.Cells.Item(5,18).Value = 1247.095685705562
$ws.Cells.Item(5,19).Value = 0.2108479789036393
$ws.Cells.Item(5,20).Value = 0.2108479789036393

# row 6: A=M1 D=FAPs
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Rtn4"
$ws.Cells.Item(6,3).Value = "Tnfrsf19"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 80.22922533333333
$ws.Cells.Item(6,8).Value = 240.687676
$ws.Cells.Item(6,9).Value = 0.194654964395505
$ws.Cells.Item(6,10).Value = 0.194654964395505
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.11008
$ws.Cells.Item(6,14).Value = 0.33024
$ws.Cells.Item(6,15).Value = 0.06903792193848049
$ws.Cells.Item(6,16).Value = 0.06903792193848049
$ws.Cells.Item(6,17).Value = 8.831633124693333
$ws.Cells.Item(6,18).Value = 79.48469812223999
$ws.Cells.Item(6,19).Value = 0.01343857423687457
$ws.Cells.Item(6,20).Value = 0.01343857423687458

# row 7: A=M1 D=sCs
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Rtn4"
$ws.Cells.Item(7,3).Value = "Tnfrsf19"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 80.22922533333333
$ws.Cells.Item(7,8).Value = 240.687676
$ws.Cells.Item(7,9).Value = 0.194654964395505
$ws.Cells.Item(7,10).Value = 0.194654964395505
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.484406
$ws.Cells.Item(7,14).Value = 4.453218
$ws.Cells.Item(7,15).Value = 0.9309620780615195
$ws.Cells.Item(7,16).Value = 0.9309620780615195
$ws.Cells.Item(7,17).Value = 119.092743460152
$ws.Cells.Item(7,18).Value = 1071.834691141368
$ws.Cells.Item(7,19).Value = 0.1812163901586304
$ws.Cells.Item(7,20).Value = 0.1812163901586305

# row 8: A=M2 D=FAPs
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Rtn4"
$ws.Cells.Item(8,3).Value = "Tnfrsf19"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 80.97090133333333
$ws.Cells.Item(8,8).Value = 242.912704
$ws.Cells.Item(8,9).Value = 0.1964544447565975
$ws.Cells.Item(8,10).Value = 0.1964544447565976
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.11008
$ws.Cells.Item(8,14).Value = 0.33024
$ws.Cells.Item(8,15).Value = 0.06903792193848049
$ws.Cells.Item(8,16).Value = 0.06903792193848049
$ws.Cells.Item(8,17).Value = 8.913276818773333
$ws.Cells.Item(8,18).Value = 80.21949136895999
$ws.Cells.Item(8,19).Value = 0.01356280662157351
$ws.Cells.Item(8,20).Value = 0.01356280662157351

# row 9: A=M2 D=sCs
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Rtn4"
$ws.Cells.Item(9,3).Value = "Tnfrsf19"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 80.97090133333333
$ws.Cells.Item(9,8).Value = 242.912704
$ws.Cells.Item(9,9).Value = 0.1964544447565975
$ws.Cells.Item(9,10).Value = 0.1964544447565976
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.484406
$ws.Cells.Item(9,14).Value = 4.453218
$ws.Cells.Item(9,15).Value = 0.9309620780615195
$ws.Cells.Item(9,16).Value = 0.9309620780615195
$ws.Cells.Item(9,17).Value = 120.193691764608
$ws.Cells.Item(9,18).Value = 1081.743225881472
$ws.Cells.Item(9,19).Value = 0.182891638135024
$ws.Cells.Item(9,20).Value = 0.1828916381350241

# row 10: A=Neutro D=FAPs
$ws.Cells.Item(10,1).Value = "Neutro"
$ws.Cells.Item(10,2).Value = "Rtn4"
$ws.Cells.Item(10,3).Value = "Tnfrsf19"
$ws.Cells.Item(10,4).Value = "FAPs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 56.835931
$ws.Cells.Item(10,8).Value = 170.507793
$ws.Cells.Item(10,9).Value = 0.1378973320410935
$ws.Cells.Item(10,10).Value = 0.1378973320410936
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.11008
$ws.Cells.Item(10,14).Value = 0.33024
$ws.Cells.Item(10,15).Value = 0.06903792193848049
$ws.Cells.Item(10,16).Value = 0.06903792193848049
$ws.Cells.Item(10,17).Value = 6.256499284479999
$ws.Cells.Item(10,18).Value = 56.30849356032
$ws.Cells.Item(10,19).Value = 0.00952014524497774
$ws.Cells.Item(10,20).Value = 0.009520145244977742

# row 11: A=Neutro D=sCs
$ws.Cells.Item(11,1).Value = "Neutro"
$ws.Cells.Item(11,2).Value = "Rtn4"
$ws.Cells.Item(11,3).Value = "Tnfrsf19"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 56.835931
$ws.Cells.Item(11,8).Value = 170.507793
$ws.Cells.Item(11,9).Value = 0.1378973320410935
$ws.Cells.Item(11,10).Value = 0.1378973320410936
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.484406
$ws.Cells.Item(11,14).Value = 4.453218
$ws.Cells.Item(11,15).Value = 0.9309620780615195
$ws.Cells.Item(11,16).Value = 0.9309620780615195
$ws.Cells.Item(11,17).Value = 84.36759699198599
$ws.Cells.Item(11,18).Value = 759.308372927874
$ws.Cells.Item(11,19).Value = 0.1283771867961158
$ws.Cells.Item(11,20).Value = 0.1283771867961158

# row 12: A=sCs D=FAPs
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Rtn4"
$ws.Cells.Item(12,3).Value = "Tnfrsf19"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 77.86108299999999
$ws.Cells.Item(12,8).Value = 233.583249
$ws.Cells.Item(12,9).Value = 0.1889092942900881
$ws.Cells.Item(12,10).Value = 0.1889092942900882
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.11008
$ws.Cells.Item(12,14).Value = 0.33024
$ws.Cells.Item(12,15).Value = 0.06903792193848049
$ws.Cells.Item(12,16).Value = 0.06903792193848049
$ws.Cells.Item(12,17).Value = 8.570948016639999
$ws.Cells.Item(12,18).Value = 77.13853214976
$ws.Cells.Item(12,19).Value = 0.01304190511265254
$ws.Cells.Item(12,20).Value = 0.01304190511265255

# row 13: A=sCs D=sCs
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Rtn4"
$ws.Cells.Item(13,3).Value = "Tnfrsf19"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 77.86108299999999
$ws.Cells.Item(13,8).Value = 233.583249
$ws.Cells.Item(13,9).Value = 0.1889092942900881
$ws.Cells.Item(13,10).Value = 0.1889092942900882
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.484406
$ws.Cells.Item(13,14).Value = 4.453218
$ws.Cells.Item(13,15).Value = 0.9309620780615195
$ws.Cells.Item(13,16).Value = 0.9309620780615195
$ws.Cells.Item(13,17).Value = 115.577458771698
$ws.Cells.Item(13,18).Value = 1040.197128945282
$ws.Cells.Item(13,19).Value = 0.1758673891774356
$ws.Cells.Item(13,20).Value = 0.1758673891774356
